$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 138 (shifts existing rows 138-154 down to 140-156)
$ws.Rows.Item(138).Insert()
$ws.Rows.Item(138).Insert()

# Populate the two new rows with the new weekly record (Primera / Segunda)
$ws.Cells.Item(138, 1).Value = 11
$ws.Cells.Item(138, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(138, 3).Value = "Bíobío"
$ws.Cells.Item(138, 4).Value = 45265
$ws.Cells.Item(138, 5).Value = 8
$ws.Cells.Item(138, 6).Value = "Fruta"
$ws.Cells.Item(138, 7).Value = 100101
$ws.Cells.Item(138, 8).Value = "Berries"
$ws.Cells.Item(138, 9).Value = 100101001
$ws.Cells.Item(138, 10).Value = "Arándano (blue)"
$ws.Cells.Item(138, 11).Value = "Sin especificar"
$ws.Cells.Item(138, 12).Value = "Primera"
$ws.Cells.Item(138, 13).Value = 100
$ws.Cells.Item(138, 14).Value = 6000
$ws.Cells.Item(138, 15).Value = 6000
$ws.Cells.Item(138, 16).Value = 6000
$ws.Cells.Item(138, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(138, 18).Value = "Región de Ñuble"
$ws.Cells.Item(138, 19).Value = 3000
$ws.Cells.Item(138, 20).Value = 2

$ws.Cells.Item(139, 1).Value = 11
$ws.Cells.Item(139, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(139, 3).Value = "Bíobío"
$ws.Cells.Item(139, 4).Value = 45265
$ws.Cells.Item(139, 5).Value = 8
$ws.Cells.Item(139, 6).Value = "Fruta"
$ws.Cells.Item(139, 7).Value = 100101
$ws.Cells.Item(139, 8).Value = "Berries"
$ws.Cells.Item(139, 9).Value = 100101001
$ws.Cells.Item(139, 10).Value = "Arándano (blue)"
$ws.Cells.Item(139, 11).Value = "Sin especificar"
$ws.Cells.Item(139, 12).Value = "Segunda"
$ws.Cells.Item(139, 13).Value = 80
$ws.Cells.Item(139, 14).Value = 5000
$ws.Cells.Item(139, 15).Value = 5000
$ws.Cells.Item(139, 16).Value = 5000
$ws.Cells.Item(139, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(139, 18).Value = "Región de Ñuble"
$ws.Cells.Item(139, 19).Value = 2500
$ws.Cells.Item(139, 20).Value = 2

